$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 14: "character" part of speech changes from "n.single" to "n."
$ws.Range("C14").Value = "n."

# Row 20: fill in "senator" entry
$ws.Range("A20").Value = "senator"
$ws.Range("B20").Value = "/'senEtE/"
$ws.Range("C20").Value = "n.C."

# Row 21: fill in "senate" entry
$ws.Range("A21").Value = "senate"
$ws.Range("B21").Value = "/'senEtE/"
$ws.Range("C21").Value = "n.Singular"

# Update selection to match the author's final cursor position
$ws.Range("C14").Select()
